# Apply updated crypto price/volume figures to Sheet1 (columns D and E).
# Values that look like plain decimal numbers (single "." separator) are
# written with a leading apostrophe so Excel stores them as TEXT, matching
# the workbook's existing inline-string cells instead of converting them to
# numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.546.96'
$ws.Range("E2").Value = '  +5.19%  '

$ws.Range("D3").Value = '1.725.91'
$ws.Range("E3").Value = '  +4.06%  '

$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = '''226.32'
$ws.Range("E5").Value = '  +3.37%  '

$ws.Range("D6").Value = '''0.5388'
$ws.Range("E6").Value = '  +2.87%  '

$ws.Range("D7").Value = '''1.004'
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '''0.2698'
$ws.Range("E8").Value = '  +0.99%  '

$ws.Range("D9").Value = '''0.06622'
$ws.Range("E9").Value = '  +4.06%  '

$ws.Range("D10").Value = '''21.72'
$ws.Range("E10").Value = '  +5.12%  '

$ws.Range("D11").Value = '''0.07768'
$ws.Range("E11").Value = '  +0.86%  '

$ws.Range("D12").Value = '''4.659'
$ws.Range("E12").Value = '  +0.05%  '

$ws.Range("D13").Value = '1.736.13'
$ws.Range("E13").Value = '  +1.74%  '

$ws.Range("D14").Value = '1.961.19'
$ws.Range("E14").Value = '  +3.94%  '

$ws.Range("E15").Value = '  +4.61%  '

$ws.Range("D16").Value = '0.0₅8302'
$ws.Range("E16").Value = '  +0.83%  '

$ws.Range("D17").Value = '''68.19'
$ws.Range("E17").Value = '  +3.79%  '

$ws.Range("D18").Value = '27.551.72'
$ws.Range("E18").Value = '  +5.25%  '

$ws.Range("D19").Value = '''225.25'
$ws.Range("E19").Value = '  +16.99%  '

$ws.Range("E20").Value = '  +0.09%  '

$ws.Range("D21").Value = '''4.758'

$ws.Range("E22").Value = '  +2.22%  '

$ws.Range("D23").Value = '''6.125'
$ws.Range("E23").Value = '  +2.38%  '

$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").Value = '''148.05'
$ws.Range("E25").Value = '  +1.51%  '

$ws.Range("D26").Value = '''1.696'
$ws.Range("E26").Value = '  +11.21%  '

$ws.Range("E27").Value = '  +3.03%  '

$ws.Range("D28").Value = '''7.432'
$ws.Range("E28").Value = '  +1.64%  '

$ws.Range("D29").Value = '''16.83'
$ws.Range("E29").Value = '  +5.14%  '

$ws.Range("D30").Value = '''0.05588'
$ws.Range("E30").Value = '  +0.67%  '

$ws.Range("E31").Value = '  +2.51%  '

$ws.Range("D32").Value = '''3.589'
$ws.Range("E32").Value = '  +3.07%  '

$ws.Range("D33").Value = '''3.474'
$ws.Range("E33").Value = '  +2.42%  '

$ws.Range("D34").Value = '''1.669'
$ws.Range("E34").Value = '  +6.26%  '

$ws.Range("D35").Value = '''0.9667'
$ws.Range("E35").Value = '  +1.29%  '

$ws.Range("D36").Value = '''2.446'
$ws.Range("E36").Value = '  +1.74%  '

$ws.Range("D37").Value = '''2.817'
$ws.Range("E37").Value = '  +1.28%  '

$ws.Range("D38").Value = '''0.5961'
$ws.Range("E38").Value = '  +4.03%  '

$ws.Range("D39").Value = '''0.01654'
$ws.Range("E39").Value = '  +3.52%  '

$ws.Range("D40").Value = '''5.902'
$ws.Range("E40").Value = '  -0.37%  '

$ws.Range("D41").Value = '''0.8607'
$ws.Range("E41").Value = '  +3.01%  '

$ws.Range("D42").Value = '1.060.84'
$ws.Range("E42").Value = '  +2.47%  '

$ws.Range("D44").Value = '''101.66'
$ws.Range("E44").Value = '  +0.48%  '

$ws.Range("D45").Value = '1.867.69'
$ws.Range("E45").Value = '  +3.91%  '

$ws.Range("E46").Value = '  +12.55%  '

$ws.Range("D47").Value = '''59.16'
$ws.Range("E47").Value = '  +1.11%  '

$ws.Range("D48").Value = '''8.208'
$ws.Range("E48").Value = '  +2.12%  '

$ws.Range("D49").Value = '''0.4431'
$ws.Range("E49").Value = '  +1.93%  '

$ws.Range("E50").Value = '  +0.41%  '

$ws.Range("D51").Value = '''0.05292'
$ws.Range("E51").Value = '  +1.02%  '
